$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "C"
$ws.Range("F4").Value = 300
$ws.Range("D10").Select()
